# Added homework 27 and 28 Feb
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing last row (167), continuing the
# daily TTV (03020) series for 2020-02-27 and 2020-02-28.
$newRows = @(
    @{ Row = 168; Timestamp = 1582761600; Date = "2020-02-27"; Id = "03020"; Name = "TTV"; Open = 0.205; High = 0.205; Low = 0.205; Close = 0.205; Vol = "-" },
    @{ Row = 169; Timestamp = 1582848000; Date = "2020-02-28"; Id = "03020"; Name = "TTV"; Open = 0.205; High = 0.205; Low = 0.205; Close = 0.205; Vol = "-" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Timestamp

    # Columns B (date) and C (id) hold text that looks numeric/date-like
    # ("2020-02-27", "03020") -- force Text format first so Excel doesn't
    # auto-convert them to a date serial / a number with its leading zero
    # stripped, then clear the formatting again so the cell keeps plain
    # (default) styling, matching the rest of the data rows.
    $ws.Cells.Item($r.Row, 2).NumberFormat = "@"
    $ws.Cells.Item($r.Row, 2).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).ClearFormats()

    $ws.Cells.Item($r.Row, 3).NumberFormat = "@"
    $ws.Cells.Item($r.Row, 3).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).ClearFormats()

    $ws.Cells.Item($r.Row, 4).Value = $r.Name
    $ws.Cells.Item($r.Row, 5).Value = $r.Open
    $ws.Cells.Item($r.Row, 6).Value = $r.High
    $ws.Cells.Item($r.Row, 7).Value = $r.Low
    $ws.Cells.Item($r.Row, 8).Value = $r.Close
    $ws.Cells.Item($r.Row, 9).Value = $r.Vol
}
